$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: update D3 description, and fill in E3/F3/G3 examples
$ws.Range("D3").Value = "Detailed breakdown of revenue streams including net interest income, non-interest income, fee income, trading revenue, wealth management income, and investment banking fees (advisory, equity underwriting, debt underwriting)"
$ws.Range("E3").Value = "Advisory revenue and M&A deal pipeline"
$ws.Range("F3").Value = "Equity and debt underwriting fees and market share"
$ws.Range("G3").Value = "FICC and equities trading revenue"

# Row 9: update D9 description, and fill in E9/F9/G9 examples
$ws.Range("D9").Value = "Performance by business segment — segment names vary by bank (e.g., CIB, CCB, Institutional Securities, Global Markets). Include revenue, earnings, and growth metrics for each reported division and sub-segment activities such as advisory, underwriting, trading, and asset management"
$ws.Range("E9").Value = "JPM: Consumer & Community Banking (CCB), Corporate & Investment Bank (CIB), Asset & Wealth Management (AWM)"
$ws.Range("F9").Value = "GS: Global Banking & Markets, Asset & Wealth Management, Platform Solutions"
$ws.Range("G9").Value = "MS: Institutional Securities, Wealth Management, Investment Management"

# Row 10: update D10 description, and fill in E10/F10/G10 examples
$ws.Range("D10").Value = "Commentary on economic environment, Federal Reserve policy and rate outlook, inflation impacts, credit cycle dynamics, consumer spending trends, and broader market conditions"
$ws.Range("E10").Value = "Federal Reserve rate trajectory and monetary policy impact"
$ws.Range("F10").Value = "Consumer sentiment, spending trends, and credit cycle positioning"
$ws.Range("G10").Value = "Capital markets environment including IPO and M&A activity levels"

# Rows that stay empty in E/F/G: clear the now-obsolete placeholder empty inline-string cells
# so they are removed from the saved worksheet XML entirely.
$emptyRows = @(2, 4, 5, 6, 7, 8, 11, 12, 13, 14, 15, 16)
foreach ($r in $emptyRows) {
    $rangeAddress = "E" + $r + ":G" + $r
    $ws.Range($rangeAddress).ClearContents()
}
